$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates ---
$ws.Range("A8").Value = 'Volume 31   Number  40'
$ws.Range("C9").Value = 'Report Covering the Week  9/30/2024  Through  10/6/2024'

# --- Cells transitioning from placeholder text style to numeric style (14 -> 15/16) ---
$ws.Range("D14").NumberFormat = '#,##0'
$ws.Range("D14").Value = 1
$ws.Range("E14").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("E14").Value = -100
$ws.Range("G14").NumberFormat = '#,##0'
$ws.Range("G14").Value = 1
$ws.Range("H14").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("H14").Value = -100
$ws.Range("D22").NumberFormat = '#,##0'
$ws.Range("D22").Value = 3
$ws.Range("E22").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("E22").Value = -100

# --- Cells transitioning from numeric style back to placeholder text style (15/16 -> 14) ---
$ws.Range("C14").Copy()
$ws.Range("D15").PasteSpecial(-4122)
$ws.Range("D15").Value = '0'
$ws.Range("C14").Copy()
$ws.Range("E15").PasteSpecial(-4122)
$ws.Range("E15").Value = '***.*'
$ws.Range("C14").Copy()
$ws.Range("D29").PasteSpecial(-4122)
$ws.Range("D29").Value = '0'
$ws.Range("C14").Copy()
$ws.Range("E29").PasteSpecial(-4122)
$ws.Range("E29").Value = '***.*'
$ws.Range("C14").Copy()
$ws.Range("D30").PasteSpecial(-4122)
$ws.Range("D30").Value = '0'
$ws.Range("C14").Copy()
$ws.Range("E30").PasteSpecial(-4122)
$ws.Range("E30").Value = '***.*'
$excel.CutCopyMode = $false

# --- Remaining simple value updates ---
$ws.Range("J14").Value = 3
$ws.Range("K14").Value = -33.333333333333
$ws.Range("F15").Value = 4
$ws.Range("G15").Value = 4
$ws.Range("H15").Value = 0
$ws.Range("I15").Value = 32
$ws.Range("K15").Value = 23.076923076923
$ws.Range("L15").Value = 14.285714285714
$ws.Range("M15").Value = 77.777777777777
$ws.Range("N15").Value = 28
$ws.Range("C16").Value = 11
$ws.Range("D16").Value = 18
$ws.Range("E16").Value = -38.888888888888
$ws.Range("F16").Value = 46
$ws.Range("G16").Value = 51
$ws.Range("H16").Value = -9.803921568627
$ws.Range("I16").Value = 406
$ws.Range("J16").Value = 362
$ws.Range("K16").Value = 12.154696132596
$ws.Range("L16").Value = 32.247557003257
$ws.Range("M16").Value = 47.101449275362
$ws.Range("N16").Value = -65.939597315436
$ws.Range("C17").Value = 17
$ws.Range("D17").Value = 20
$ws.Range("E17").Value = -15
$ws.Range("F17").Value = 55
$ws.Range("G17").Value = 66
$ws.Range("H17").Value = -16.666666666666
$ws.Range("I17").Value = 637
$ws.Range("J17").Value = 573
$ws.Range("K17").Value = 11.169284467713
$ws.Range("L17").Value = 55.365853658536
$ws.Range("M17").Value = 231.770833333333
$ws.Range("N17").Value = 68.518518518518
$ws.Range("C18").Value = 4
$ws.Range("D18").Value = 4
$ws.Range("E18").Value = 0
$ws.Range("F18").Value = 17
$ws.Range("G18").Value = 14
$ws.Range("H18").Value = 21.428571428571
$ws.Range("I18").Value = 189
$ws.Range("J18").Value = 149
$ws.Range("K18").Value = 26.845637583892
$ws.Range("L18").Value = 45.384615384615
$ws.Range("M18").Value = -13.698630136986
$ws.Range("N18").Value = -88.842975206611
$ws.Range("C19").Value = 29
$ws.Range("D19").Value = 31
$ws.Range("E19").Value = -6.451612903225
$ws.Range("F19").Value = 97
$ws.Range("G19").Value = 84
$ws.Range("H19").Value = 15.47619047619
$ws.Range("I19").Value = 870
$ws.Range("J19").Value = 824
$ws.Range("K19").Value = 5.582524271844
$ws.Range("L19").Value = -5.639913232104
$ws.Range("M19").Value = 110.653753026634
$ws.Range("N19").Value = -7.545164718384
$ws.Range("C20").Value = 2
$ws.Range("D20").Value = 13
$ws.Range("E20").Value = -84.615384615384
$ws.Range("F20").Value = 13
$ws.Range("G20").Value = 41
$ws.Range("H20").Value = -68.292682926829
$ws.Range("I20").Value = 224
$ws.Range("J20").Value = 267
$ws.Range("K20").Value = -16.104868913857
$ws.Range("L20").Value = 33.333333333333
$ws.Range("M20").Value = 96.491228070175
$ws.Range("N20").Value = -87.081891580161
$ws.Range("C21").Value = 65
$ws.Range("D21").Value = 87
$ws.Range("E21").Value = -25.287356321839
$ws.Range("F21").Value = 232
$ws.Range("G21").Value = 261
$ws.Range("H21").Value = -11.111111111111
$ws.Range("I21").Value = 2360
$ws.Range("J21").Value = 2204
$ws.Range("K21").Value = 7.078039927404
$ws.Range("L21").Value = 19.79695431472
$ws.Range("M21").Value = 91.247974068071
$ws.Range("N21").Value = -60.581259395356
$ws.Range("F22").Value = 2
$ws.Range("G22").Value = 5
$ws.Range("H22").Value = -60
$ws.Range("J22").Value = 40
$ws.Range("K22").Value = -2.5
$ws.Range("L22").Value = 30
$ws.Range("D24").Value = 45
$ws.Range("E24").Value = -4.444444444444
$ws.Range("F24").Value = 199
$ws.Range("G24").Value = 209
$ws.Range("H24").Value = -4.784688995215
$ws.Range("I24").Value = 2382
$ws.Range("J24").Value = 2158
$ws.Range("K24").Value = 10.379981464318
$ws.Range("L24").Value = 31.456953642384
$ws.Range("M24").Value = 77.893950709484
$ws.Range("C25").Value = 23
$ws.Range("D25").Value = 34
$ws.Range("E25").Value = -32.35294117647
$ws.Range("F25").Value = 144
$ws.Range("G25").Value = 155
$ws.Range("H25").Value = -7.096774193548
$ws.Range("I25").Value = 1870
$ws.Range("J25").Value = 1577
$ws.Range("K25").Value = 18.57958148383
$ws.Range("L25").Value = 45.525291828793
$ws.Range("C26").Value = 19
$ws.Range("D26").Value = 24
$ws.Range("E26").Value = -20.833333333333
$ws.Range("F26").Value = 107
$ws.Range("G26").Value = 99
$ws.Range("H26").Value = 8.080808080808
$ws.Range("I26").Value = 1053
$ws.Range("J26").Value = 896
$ws.Range("K26").Value = 17.522321428571
$ws.Range("L26").Value = 52.608695652173
$ws.Range("M26").Value = 113.157894736842
$ws.Range("C27").Value = 2
$ws.Range("E27").Value = 100
$ws.Range("F27").Value = 6
$ws.Range("G27").Value = 6
$ws.Range("H27").Value = 0
$ws.Range("I27").Value = 47
$ws.Range("J27").Value = 41
$ws.Range("K27").Value = 14.634146341463
$ws.Range("L27").Value = 2.173913043478
$ws.Range("I28").Value = 119
$ws.Range("J28").Value = 108
$ws.Range("K28").Value = 10.185185185185
$ws.Range("L28").Value = 26.595744680851
$ws.Range("N29").Value = -95.918367346938
$ws.Range("N30").Value = -95.555555555555
